$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.730.86'
$ws.Range("E2").Value = '  -4.31%  '
$ws.Range("D3").Value = '3.273.20'
$ws.Range("E3").Value = '  -4.74%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '555.87'
$ws.Range("E5").Value = '  -2.78%  '
$ws.Range("D6").Value = '183.45'
$ws.Range("E6").Value = '  -2.98%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").Value = '  -2.51%  '
$ws.Range("D9").Value = '3.271.06'
$ws.Range("E9").Value = '  -4.53%  '
$ws.Range("D10").Value = '0.187'
$ws.Range("E10").Value = '  -6.94%  '
$ws.Range("D11").Value = '0.587'
$ws.Range("E11").Value = '  -3.75%  '
$ws.Range("D12").Value = '47.30'
$ws.Range("E12").Value = '  -6.76%  '
$ws.Range("D13").Value = '0.0000266'
$ws.Range("E13").Value = '  -5.46%  '
$ws.Range("D14").Value = '633.75'
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("D15").Value = '8.57'
$ws.Range("E15").Value = '  -4.53%  '
$ws.Range("D16").Value = '3.796.49'
$ws.Range("E16").Value = '  -4.47%  '
$ws.Range("D17").Value = '65.689.81'
$ws.Range("E17").Value = '  -4.15%  '
$ws.Range("D18").Value = '17.85'
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("E19").Value = '  -3.19%  '
$ws.Range("D20").Value = '3.266.21'
$ws.Range("E20").Value = '  -4.87%  '
$ws.Range("D21").Value = '11.37'
$ws.Range("E21").Value = '  -6.21%  '
$ws.Range("D22").Value = '0.904'
$ws.Range("E22").Value = '  -2.99%  '
$ws.Range("D23").Value = '17.64'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '106.21'
$ws.Range("E24").Value = '  +8.46%  '
$ws.Range("D25").Value = '4.92'
$ws.Range("E25").Value = '  -6.35%  '
$ws.Range("D26").Value = '3.97'
$ws.Range("E26").Value = '  -5.77%  '
$ws.Range("D27").Value = '2.66'
$ws.Range("E27").Value = '  -5.40%  '
$ws.Range("D28").Value = '9.56'
$ws.Range("E28").Value = '  -2.05%  '
$ws.Range("D29").Value = '8.67'
$ws.Range("E29").Value = '  -4.79%  '
$ws.Range("D30").Value = '30.41'
$ws.Range("E30").Value = '  -5.10%  '
$ws.Range("D31").Value = '4.05'
$ws.Range("E31").Value = '  -1.38%  '
$ws.Range("D32").Value = '6.30'
$ws.Range("E32").Value = '  -4.77%  '
$ws.Range("D33").Value = '11.00'
$ws.Range("E33").Value = '  -4.05%  '
$ws.Range("D34").Value = '544.00'
$ws.Range("E34").Value = '  +13.10%  '
$ws.Range("D35").Value = '0.105'
$ws.Range("E35").Value = '  -2.70%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '57.26'
$ws.Range("E37").Value = '  -5.68%  '
$ws.Range("D38").Value = '3.688.07'
$ws.Range("E38").Value = '  +1.57%  '
$ws.Range("D39").Value = '3.39'
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0722'
$ws.Range("E40").Value = '  -6.45%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.130'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("D42").Value = '2.72'
$ws.Range("E42").Value = '  -4.25%  '
$ws.Range("D43").Value = '3.35'
$ws.Range("E43").Value = '  -2.98%  '
$ws.Range("D44").Value = '32.39'
$ws.Range("E44").Value = '  -4.40%  '
$ws.Range("D45").Value = '0.337'
$ws.Range("E45").Value = '  -7.37%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0414'
$ws.Range("E46").Value = '  -4.58%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.24'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("E48").Value = '  -5.56%  '
$ws.Range("D49").Value = '0.129'
$ws.Range("E49").Value = '  -3.00%  '
$ws.Range("D50").Value = '0.999'
$ws.Range("D51").Value = '1.24'
$ws.Range("E51").Value = '  +2.36%  '
